# Sync attendance_reports: fix "Recorded By" column (G) ordering of names.
# Replace every exact occurrence of "System, dnasr281@gmail.com" with
# "dnasr281@gmail.com, System" in column G of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colG = $ws.Range("G1:G$lastRow")

for ($i = 1; $i -le $colG.Cells.Count; $i++) {
    $cell = $colG.Cells.Item($i)
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value2 = $newVal
    }
}
